$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESCALETA")
$ws.Columns("N").Cut() | Out-Null
$ws.Columns("M").Insert() | Out-Null

# The merged header label ("Tipo de Motor") anchored at M1 travelled with the
# physical column move; put it back on the top-left cell of the M1:N1 merge.
$label = $ws.Range("N1").Value
$ws.Range("N1").Value = $null
$ws.Range("M1").Value = $label

# Re-apply a box border around the merged M1:N1 header cell so the internal
# divider between M1 and N1 is removed (matches Excel's behaviour when you
# draw an outside border on a merged range).
$headerRange = $ws.Range("M1:N1")
$headerRange.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$headerRange.Borders.Item(7).Weight = -4138  # xlMedium
$headerRange.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$headerRange.Borders.Item(10).Weight = -4138
$headerRange.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$headerRange.Borders.Item(8).Weight = -4138
$headerRange.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$headerRange.Borders.Item(9).Weight = -4138
$headerRange.Borders.Item(11).LineStyle = 0  # xlInsideVertical - none
$headerRange.Borders.Item(12).LineStyle = 0  # xlInsideHorizontal - none

$ws.Range("M1:N1").Select() | Out-Null
